# Updating projet_list_all with project test set (3 projects per instrument)
#
# Adds a new "object_annotation_category" value in column Y (the
# Category.field column) for every data row of the "Data" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 47
}

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Range("Y$r").Value = "object_annotation_category"
}
